$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Datos crudos")
$ws2 = $wb.Worksheets.Item("Datos válidos")

# Update Timestamp (C) and Temperature (E) columns for rows 2-33
$ws1.Range("C2").Value = "2023-12-09 07:16:32"
$ws1.Range("E2").Value2 = 25.673076923076898
$ws1.Range("C3").Value = "2023-12-09 07:17:33"
$ws1.Range("E3").Value2 = 23.837412587412501
$ws1.Range("C4").Value = "2023-12-09 07:18:34"
$ws1.Range("E4").Value2 = 22.788461538461501
$ws1.Range("C5").Value = "2023-12-09 07:19:35"
$ws1.Range("E5").Value2 = 21.870629370629299
$ws1.Range("C6").Value = "2023-12-09 07:20:37"
$ws1.Range("E6").Value2 = 21.739510489510501
$ws1.Range("C7").Value = "2023-12-09 07:21:38"
$ws1.Range("E7").Value2 = 21.346153846153801
$ws1.Range("C8").Value = "2023-12-09 07:22:39"
$ws1.Range("E8").Value2 = 21.083916083916101
$ws1.Range("C9").Value = "2023-12-09 07:23:41"
$ws1.Range("E9").Value2 = 20.5594405594405
$ws1.Range("C10").Value = "2023-12-09 07:24:42"
$ws1.Range("E10").Value2 = 20.5594405594405
$ws1.Range("C11").Value = "2023-12-09 07:25:43"
$ws1.Range("E11").Value2 = 20.2972027972028
$ws1.Range("C12").Value = "2023-12-09 07:26:45"
$ws1.Range("E12").Value2 = 19.9038461538461
$ws1.Range("C13").Value = "2023-12-09 07:27:46"
$ws1.Range("E13").Value2 = 20.034965034965001
$ws1.Range("C14").Value = "2023-12-09 07:28:47"
$ws1.Range("E14").Value2 = 19.379370629370602
$ws1.Range("C15").Value = "2023-12-09 07:29:49"
$ws1.Range("E15").Value2 = 19.510489510489499
$ws1.Range("C16").Value = "2023-12-09 07:30:50"
$ws1.Range("E16").Value2 = 19.510489510489499
$ws1.Range("C17").Value = "2023-12-09 07:31:51"
$ws1.Range("E17").Value2 = 19.9038461538461
$ws1.Range("C18").Value = "2023-12-09 07:32:53"
$ws1.Range("E18").Value2 = 19.2482517482517
$ws1.Range("C19").Value = "2023-12-09 07:33:54"
$ws1.Range("E19").Value2 = 19.510489510489499
$ws1.Range("C20").Value = "2023-12-09 07:34:55"
$ws1.Range("E20").Value2 = 19.379370629370602
$ws1.Range("C21").Value = "2023-12-09 07:35:57"
$ws1.Range("E21").Value2 = 19.9038461538461
$ws1.Range("C22").Value = "2023-12-09 07:36:58"
$ws1.Range("E22").Value2 = 19.379370629370602
$ws1.Range("C23").Value = "2023-12-09 07:37:59"
$ws1.Range("E23").Value2 = 19.641608391608401
$ws1.Range("C24").Value = "2023-12-09 07:39:01"
$ws1.Range("E24").Value2 = 19.510489510489499
$ws1.Range("C25").Value = "2023-12-09 07:40:02"
$ws1.Range("E25").Value2 = 19.641608391608401
$ws1.Range("C26").Value = "2023-12-09 07:41:03"
$ws1.Range("E26").Value2 = 19.510489510489499
$ws1.Range("C27").Value = "2023-12-09 07:42:05"
$ws1.Range("E27").Value2 = 19.379370629370602
$ws1.Range("C28").Value = "2023-12-09 07:43:06"
$ws1.Range("E28").Value2 = 19.510489510489499
$ws1.Range("C29").Value = "2023-12-09 07:44:07"
$ws1.Range("E29").Value2 = 19.641608391608401
$ws1.Range("C30").Value = "2023-12-09 07:45:09"
$ws1.Range("E30").Value2 = 19.379370629370602
$ws1.Range("C31").Value = "2023-12-09 07:46:10"
$ws1.Range("E31").Value2 = 19.379370629370602
$ws1.Range("C32").Value = "2023-12-09 07:47:11"
$ws1.Range("E32").Value2 = 19.510489510489499
$ws1.Range("C33").Value = "2023-12-09 07:48:13"
$ws1.Range("E33").Value2 = 20.034965034965001

# Remove the now-obsolete tail rows (old rows 34-45)
$ws1.Rows("34:45").Delete()

# Update settling-time lookup row index and valid-row counter formula
$ws1.Range("H2").Value2 = 14
$ws1.Range("H3").Formula = "=COUNT(E:E)-1"

# Restore view/selection state (sheet1 ends up active/selected, matching the workbook)
$ws2.Activate()
$ws2.Range("B3").Select()
$ws1.Activate()
$ws1.Range("E19").Select()
